$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new day's price record for "Región de Arica y Parinacota" was added to the
# weekly sheet. It becomes the new row 291, pushing every existing row from
# 291 downward down by one (291 -> 292, ..., 369 -> 370).
$ws.Rows.Item(291).Insert()

$ws.Range("A291").Value = 3
$ws.Range("B291").Value = "Femacal de La Calera"
$ws.Range("C291").Value = "Coquimbo"
$ws.Range("D291").Value = 44736
$ws.Range("E291").Value = 5
$ws.Range("F291").Value = 100112043
$ws.Range("G291").Value = "Pepino ensalada"
$ws.Range("H291").Value = "Sin especificar"
$ws.Range("I291").Value = "Primera"
$ws.Range("J291").Value = 85
$ws.Range("K291").Value = 18000
$ws.Range("L291").Value = 19000
$ws.Range("M291").Value = 18529
$ws.Range("N291").Value = "$/caja 70 unidades"
$ws.Range("O291").Value = "Región de Arica y Parinacota"
$ws.Range("P291").Value = 265
$ws.Range("Q291").Value = 70
$ws.Range("R291").Value = "Hortaliza"
